$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header title strings (rich text runs with identical formatting; plain text is equivalent) ---
$ws.Range("A8").Value = "Volume 32   Number  30"
$ws.Range("C9").Value = "Report Covering the Week  7/21/2025  Through  7/27/2025"

# --- Simple numeric value updates (type/style unchanged) ---
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = -37.5
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -80
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 61
$ws.Range("J16").Value = 93
$ws.Range("K16").Value = -34.408602150537
$ws.Range("L16").Value = -43.518518518518
$ws.Range("M16").Value = -19.736842105263
$ws.Range("N16").Value = -87.103594080338
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 500
$ws.Range("F17").Value = 31
$ws.Range("H17").Value = 14.814814814814
$ws.Range("I17").Value = 83
$ws.Range("J17").Value = 95
$ws.Range("K17").Value = -12.631578947368
$ws.Range("L17").Value = -33.6
$ws.Range("M17").Value = 59.615384615384
$ws.Range("N17").Value = -49.696969696969
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 133.333333333333
$ws.Range("F18").Value = 21
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = 61.538461538461
$ws.Range("I18").Value = 97
$ws.Range("J18").Value = 148
$ws.Range("K18").Value = -34.459459459459
$ws.Range("L18").Value = -46.111111111111
$ws.Range("M18").Value = -11.009174311926
$ws.Range("N18").Value = -78.958785249457
$ws.Range("C19").Value = 26
$ws.Range("D19").Value = 20
$ws.Range("E19").Value = 30
$ws.Range("F19").Value = 89
$ws.Range("G19").Value = 115
$ws.Range("H19").Value = -22.608695652173
$ws.Range("I19").Value = 542
$ws.Range("J19").Value = 615
$ws.Range("K19").Value = -11.869918699187
$ws.Range("L19").Value = -22.460658082975
$ws.Range("M19").Value = -10.413223140495
$ws.Range("N19").Value = -60.205580029368
$ws.Range("L20").Value = -60.869565217391
$ws.Range("M20").Value = -60.869565217391
$ws.Range("N20").Value = -97.831325301204
$ws.Range("C21").Value = 40
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = 37.931034482758
$ws.Range("F21").Value = 152
$ws.Range("G21").Value = 170
$ws.Range("H21").Value = -10.588235294117
$ws.Range("I21").Value = 798
$ws.Range("J21").Value = 976
$ws.Range("K21").Value = -18.237704918032
$ws.Range("L21").Value = -30.061349693251
$ws.Range("M21").Value = -8.381171067738
$ws.Range("N21").Value = -72.330097087378
$ws.Range("H22").Value = -100
$ws.Range("J22").Value = 27
$ws.Range("K22").Value = 11.111111111111
$ws.Range("L22").Value = 15.384615384615
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 43
$ws.Range("E24").Value = -48.837209302325
$ws.Range("F24").Value = 105
$ws.Range("G24").Value = 167
$ws.Range("H24").Value = -37.125748502994
$ws.Range("I24").Value = 813
$ws.Range("J24").Value = 995
$ws.Range("K24").Value = -18.291457286432
$ws.Range("L24").Value = -31.334459459459
$ws.Range("M24").Value = -7.718501702610
$ws.Range("C25").Value = 19
$ws.Range("D25").Value = 34
$ws.Range("E25").Value = -44.117647058823
$ws.Range("F25").Value = 79
$ws.Range("G25").Value = 122
$ws.Range("H25").Value = -35.245901639344
$ws.Range("I25").Value = 595
$ws.Range("J25").Value = 805
$ws.Range("K25").Value = -26.086956521739
$ws.Range("L25").Value = -33.667781493868
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 350
$ws.Range("F26").Value = 39
$ws.Range("G26").Value = 31
$ws.Range("H26").Value = 25.806451612903
$ws.Range("I26").Value = 220
$ws.Range("J26").Value = 209
$ws.Range("K26").Value = 5.263157894736
$ws.Range("L26").Value = -14.396887159533
$ws.Range("M26").Value = 58.273381294964
$ws.Range("H27").Value = -100
$ws.Range("L27").Value = 14.285714285714
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 80
$ws.Range("I28").Value = 46
$ws.Range("J28").Value = 40
$ws.Range("K28").Value = 15
$ws.Range("L28").Value = 17.948717948717
$ws.Range("F29").Value = 1
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 3
$ws.Range("K29").Value = 200
$ws.Range("M29").Value = 200
$ws.Range("N29").Value = -57.142857142857
$ws.Range("I30").Value = 2
$ws.Range("K30").Value = 100
$ws.Range("M30").Value = 100
$ws.Range("N30").Value = -71.428571428571
$ws.Range("L31").Value = -12.5

# --- Cells changing between numeric and text ("0" / "***.*") representation ---
# Use a same-column/row-type reference cell to copy number format/style, then set the value
$ws.Range("C20").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C20").PasteSpecial(-4122)

$ws.Range("D20").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D20").PasteSpecial(-4122)

$ws.Range("E20").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E20").PasteSpecial(-4122)

$ws.Range("F22").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("F22").PasteSpecial(-4122)

$ws.Range("F27").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("F27").PasteSpecial(-4122)

$ws.Range("I14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C29").Value = 1

$ws.Range("I14").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C30").Value = 1

$ws.Range("G31").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("G31").PasteSpecial(-4122)

$ws.Range("H31").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("H31").PasteSpecial(-4122)

$excel.CutCopyMode = 0